$wb = $excel.ActiveWorkbook
$tc020 = $wb.Worksheets.Item("tc020")

# 1) Consume sheetId 41 via a throwaway sheet so the new "tc010" sheet ends up
#    with sheetId 42, matching the target workbook.xml.
$dummy = $wb.Worksheets.Add([System.Type]::Missing, $tc020)
$dummy.Name = "DummyX"

# 2) Duplicate tc020 in full (keeps its xmlns/sheetFormatPr/styles/selection
#    intact) right after the dummy sheet, i.e. at the very end of the tab list.
$tc020.Copy([System.Type]::Missing, $dummy) | Out-Null

# 3) Remove the throwaway sheet now that its sheetId has been spent.
$dummy.Delete() | Out-Null

# 4) Rename the duplicate to "tc010" and drop its defect column (F), leaving
#    just columns A:E (Projectname/ReleaseName/CycleName/SuiteName/TCClick).
$tc010 = $wb.Worksheets.Item("tc020 (2)")
$tc010.Name = "tc010"
$tc010.Columns.Item(6).Delete() | Out-Null

# 5) Fix up the defect id/summary text on tc020's F column.
$tc020.Range("F1").Value = "DF-317"
$tc020.Range("F2").Value = "defid"

# 6) Re-fetch sheet references (collection shifted after the add/copy/delete)
#    and restore the view/selection state shown in the target workbook.
$tc020 = $wb.Worksheets.Item("tc020")
$tc010 = $wb.Worksheets.Item("tc010")

$tc020.Range("A1:E2").Select() | Out-Null

$tc010.Activate() | Out-Null
$tc010.Range("F7").Select() | Out-Null
